# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price cells in column D are plain-decimal-looking text (e.g. "608.67"); Excel
# would otherwise auto-coerce those to numbers on assignment, so for those we
# force Text format first and restore the default "Normal" style afterwards
# (keeps the cell un-styled, matching the rest of the sheet, while the stored
# value stays a string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.925.06'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '3.511.97'
$ws.Range('E3').Value = '  +0.27%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '608.67'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '147.83'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('D7').Value = '3.512.68'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.78%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.143'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.06%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '8.02'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +5.81%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.423'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.85%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000218'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '31.98'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '4.105.62'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '3.508.12'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').Value = '66.937.23'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('E18').Value = '  -0.27%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.74'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +8.04%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.46'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.68%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '15.36'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.31%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '438.15'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.88%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.609'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.49%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '79.65'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').Value = '3.647.74'
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('E26').Value = '  +0.07%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.0000121'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -4.15%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '9.79'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.87%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '8.25'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -4.25%  '
$ws.Range('E31').Value = '  -2.42%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.168'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('E33').Value = '  -0.03%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '25.58'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.97'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -2.86%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.81'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.13%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '8.07'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.02%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.02%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '175.50'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.79%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0894'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('E42').Value = '  -0.20%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.06'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -11.80%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.895'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.00%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '46.10'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.16%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '28.10'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -7.03%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.26'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -2.35%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.48'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.81%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -3.07%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.994'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.248'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.67%  '
